{"js": "// Insert a new \"Project\" bullet item (AAR671 Advanced Computer Architecture...)\n// immediately before the existing \"ELE510 Advanced Digital Logic Design...\"\n// bullet item, reusing that paragraph's list/style formatting.\n\nconst anchorText =\n  \"ELE510 Advanced Digital Logic Design: Instruction Set Extension for a RISC CPU\";\n\nconst results = context.document.body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph text not found: \" + anchorText);\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\n\n// insertParagraph(\"Before\") clones the anchor paragraph's formatting\n// (style, numbering, run properties) for the new paragraph.\nanchorParagraph.insertParagraph(\n  \"AAR671 Advanced Computer Architecture: Analysis and optimization of BoomV1 architecture (Super scalar)\",\n  \"Before\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new \"Project\" bullet item (AAR671 Advanced Computer Architecture...)\n# immediately before the existing \"ELE510 Advanced Digital Logic Design...\"\n# bullet item, reusing that paragraph's list/style formatting.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"ELE510 Advanced Digital Logic Design: Instruction Set Extension for a RISC CPU\"\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Anchor paragraph text not found: $anchorText\"\n}\n\n$anchorPara = $searchRange.Paragraphs(1)\n\n# InsertParagraphBefore clones the anchor paragraph's formatting\n# (style, numbering, run properties) onto the new, empty paragraph.\n$anchorPara.Range.InsertParagraphBefore()\n\n# Re-locate the freshly inserted (still empty) paragraph: it now sits\n# immediately before the ELE510 paragraph.\n$searchRange2 = $d.Content\n$found2 = $searchRange2.Find.Execute($anchorText)\nif (-not $found2) {\n    throw \"Anchor paragraph text not found after insertion: $anchorText\"\n}\n$elePara = $searchRange2.Paragraphs(1)\n$newPara = $elePara.Previous()\n\n$newPara.Range.Text = \"AAR671 Advanced Computer Architecture: Analysis and optimization of BoomV1 architecture (Super scalar)\"\n"}
